$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "EV001"
$ws.Range("A2").Value = "自動"
$ws.Range("B2").Value = "自動"
$ws.Range("A3").Value = "ーーーーーーーーアイテム生成数ーーーーーーーー"
$ws.Range("A4").ClearContents()
